# Regenerate the handback-status report: the two source files were
# re-processed by the pipeline under new GUIDs / new xliff hashes, and
# got fresh handoff/handback timestamps.
#
# Old file 1: 1b1858c3-3b00-47a2-8f24-973e2b5359c5.md
# New file 1: 49dbc38d-2320-4cc0-81ef-5fd882349c85.md
#
# Old file 2: c418e212-b76f-4bce-8821-99186b269b0e.md
# New file 2: ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md
#
# Old xliff hash (zh-cn, row 2): 12bba9edd97355a63be4dcec563f4f1a1f89b5b0
# Old xliff hash (zh-cn, row 3): 5917ff04743a22e1bdeea04a95e466e37b6c7d73
# New xliff hash (both rows)   : 193e86ffbf271f42e1c4866a2fdea8c483367546

$wb = $excel.ActiveWorkbook

$guid1New = "49dbc38d-2320-4cc0-81ef-5fd882349c85"
$guid2New = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e"

$file1New = $guid1New + ".md"
$file2New = $guid2New + ".md"

# Helper: update a hyperlinked cell's text while keeping its Address (rId)
# and without adding a duplicate Hyperlinks entry.
function Set-HyperlinkCellText($ws, $cellRef, $newText) {
    $target = $ws.Range($cellRef).Address()
    $links = @($ws.Hyperlinks)
    foreach ($lnk in $links) {
        if ($lnk.Range.Address() -eq $target) {
            $addr = $lnk.Address
            $ws.Range($cellRef).Value = $newText
            $lnk.TextToDisplay = $newText
            $lnk.Address = $addr
            return
        }
    }
    # No hyperlink on that cell - just a plain value update.
    $ws.Range($cellRef).Value = $newText
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HyperlinkCellText $wsOverview "A2" $file1New
Set-HyperlinkCellText $wsOverview "B2" ("e2e\" + $file1New)
$wsOverview.Range("G2").Value = "2016-09-02 11:13:55"

Set-HyperlinkCellText $wsOverview "A3" $file2New
Set-HyperlinkCellText $wsOverview "B3" ("e2e\" + $file2New)
$wsOverview.Range("G3").Value = "2016-09-02 11:13:55"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = $guid1New + ".193e86ffbf271f42e1c4866a2fdea8c483367546.zh-cn.xlf"

Set-HyperlinkCellText $wsZh "A2" $file1New
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = "2016-09-02 11:13:51"
Set-HyperlinkCellText $wsZh "I2" $file1New
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = "2016-09-02 11:14:17"

Set-HyperlinkCellText $wsZh "A3" $file2New
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = "2016-09-02 11:13:51"
Set-HyperlinkCellText $wsZh "I3" $file2New
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = "2016-09-02 11:14:17"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = $guid1New + ".193e86ffbf271f42e1c4866a2fdea8c483367546.de-de.xlf"

Set-HyperlinkCellText $wsDe "A2" $file1New
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = "2016-09-02 11:13:55"
Set-HyperlinkCellText $wsDe "I2" $file1New
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = "2016-09-02 11:14:25"

Set-HyperlinkCellText $wsDe "A3" $file2New
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = "2016-09-02 11:13:55"
Set-HyperlinkCellText $wsDe "I3" $file2New
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = "2016-09-02 11:14:25"
